$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.36
$ws.Range("C2").Value = 0.34
$ws.Range("D2").Value = 0.71
$ws.Range("E2").Value = 1.2
$ws.Range("F2").Value = 1.31
$ws.Range("G2").Value = 0.61
$ws.Range("H2").Value = 26.73

$ws.Range("B3").Value = 0.44
$ws.Range("C3").Value = 0.46
$ws.Range("D3").Value = 0.73
$ws.Range("E3").Value = 1.22
$ws.Range("F3").Value = 1.04
$ws.Range("G3").Value = 0.5600000000000001
$ws.Range("H3").Value = -32.81
